# Generate Report for Handback
# Updates the handoff/handback timestamp strings in the handback-status
# workbook to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for the first (de-de) row
$wsOverview.Range("G2").Value = "2016-09-06 01:12:00"

# zh-cn!H2 - "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-06 01:11:55"

# zh-cn!K2 - "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-09-06 01:12:18"

# de-de!H2 - "Correspond Handoff Datetime" (shares text with Overview!G2)
$wsDeDe.Range("H2").Value = "2016-09-06 01:12:00"

# de-de!K2 - "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-09-06 01:12:26"
